$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "56.961.18"
$ws.Range("E2").Value = "  +0.55%  "
$ws.Range("D3").Value = "2.397.62"
$ws.Range("E3").Value = "  +0.73%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "505.97"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.91%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "133.28"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.25%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.09%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.553"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.37%  "
$ws.Range("D9").Value = "2.408.48"
$ws.Range("E9").Value = "  +0.03%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0978"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.09%  "
$ws.Range("E11").Value = "  -0.66%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.324"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.57%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.64"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.85%  "
$ws.Range("D14").Value = "2.826.20"
$ws.Range("E14").Value = "  -0.18%  "
$ws.Range("D15").Value = "56.908.78"
$ws.Range("E15").Value = "  +0.45%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "21.81"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.59%  "
$ws.Range("E17").Value = "  +2.13%  "
$ws.Range("D18").Value = "2.408.40"
$ws.Range("E18").Value = "  +0.94%  "
$ws.Range("E19").Value = "  -0.09%  "
$ws.Range("E20").Value = "  +0.03%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "310.94"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.77%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.26"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.79%  "
$ws.Range("E23").Value = "  -0.23%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.56"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -5.06%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "67.64"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.89%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.997"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.26%  "
$ws.Range("E27").Value = "  -0.35%  "
$ws.Range("E28").Value = "  -0.67%  "
$ws.Range("E29").Value = "  +2.39%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "175.96"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.66%  "
$ws.Range("E31").Value = "  +1.13%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.67"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.69%  "
$ws.Range("E33").Value = "  +0.21%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.86"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.67%  "
$ws.Range("E35").Value = "  +0.05%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.998"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.29%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "17.98"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.17%  "
$ws.Range("E38").Value = "  -1.17%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.83"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.13%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.830"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +5.19%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "36.89"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.81%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.44"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.52%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "131.79"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.15%  "
$ws.Range("E44").Value = "  +0.85%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.93"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.43%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.571"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.71%  "
$ws.Range("B47").Value = "Bittensor"
$ws.Range("C47").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "252.13"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.45%  "
$ws.Range("B48").Value = "Stellar"
$ws.Range("C48").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0915"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.45%  "
$ws.Range("E49").Value = "  -0.20%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0210"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.26%  "
$ws.Range("E51").Value = "  +7.43%  "
